$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Exclusive)
$ws.Range("B2").Value = 140.1408
$ws.Range("C2").Value = 2114.3664
$ws.Range("D2").Value = 141.162
$ws.Range("E2").Value = 69.0453
$ws.Range("F2").Value = 2464.7145

# Row 3 (Overlaps)
$ws.Range("B3").Value = 47.3189
$ws.Range("C3").Value = 1079.0288
$ws.Range("D3").Value = 741.1805
$ws.Range("E3").Value = 204.8602
$ws.Range("F3").Value = 2072.3884
